# Weekly update: add a new price record for "Poroto granado" at the top of
# the data table. This is a new row inserted above the current row 67,
# which pushes the existing rows 67-151 down to 68-152 (dimension grows
# from A1:R151 to A1:R152).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 67 (shifts rows 67..151 down to 68..152).
$ws.Rows.Item(67).Insert()

# Populate the new row 67 with this week's record.
$ws.Range("A67").Value = 5
$ws.Range("B67").Value = "Macroferia Regional de Talca"
$ws.Range("C67").Value = "Maule"
$ws.Range("D67").Value = 44915
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = 100112030
$ws.Range("G67").Value = "Poroto granado"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 150
$ws.Range("K67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("M67").Value = 35000
$ws.Range("N67").Value = "$/saco 25 kilos"
$ws.Range("O67").Value = "Región del Maule"
$ws.Range("P67").Value = 1400
$ws.Range("Q67").Value = 25
$ws.Range("R67").Value = "Hortaliza"
